# SEXTO COMMIT - Edição da documentação do jogo
#
# 1) First table, row 1, column 2: "Jogo do Unity" -> "O Capitão do Moby Dick"
#    (also drops the spell-check proofErr markers and introduces a fresh
#    _GoBack bookmark right after "O ").
# 2) Table cell containing "Não pode atacar e se movimentar ao mesmo tempo.":
#    remove the stray _GoBack bookmark that trails the text.

$d = $word.ActiveDocument
$table = $d.Tables.Item(1)

# ---- Helper: find the (row, col) of a cell whose text contains $needle ----
function Find-CellByText($tbl, $needle) {
    for ($r = 1; $r -le $tbl.Rows.Count; $r++) {
        $row = $tbl.Rows.Item($r)
        for ($c = 1; $c -le $row.Cells.Count; $c++) {
            $cell = $row.Cells.Item($c)
            if ($cell.Range.Text -like ("*" + $needle + "*")) {
                return $cell
            }
        }
    }
    return $null
}

# ---------------------------------------------------------------------
# 1) "Jogo do Unity" -> "O Capitão do Moby Dick" (+ fresh _GoBack bookmark)
# ---------------------------------------------------------------------
$titleCell = Find-CellByText $table "Jogo do"
$titlePara = $titleCell.Range.Paragraphs.Item(1)
$titleRange = $d.Range($titlePara.Range.Start, $titlePara.Range.End)

$titleXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="00EC05A3" w:rsidRDefault="00EC05A3"><w:r><w:t xml:space="preserve">O </w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:t>Capitão do Moby Dick</w:t></w:r></w:p>
'@
$titleRange.InsertXML($titleXml)

# ---------------------------------------------------------------------
# 2) Drop the leftover _GoBack bookmark after "Não pode atacar..."
# ---------------------------------------------------------------------
$restrCell = Find-CellByText $table "movimentar ao mesmo tempo"
$restrPara = $restrCell.Range.Paragraphs.Item(1)
$restrRange = $d.Range($restrPara.Range.Start, $restrPara.Range.End)

$restrXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="007C230F" w:rsidRDefault="007C230F" w:rsidP="007C230F"><w:pPr><w:pStyle w:val="PargrafodaLista"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t>Não pode atacar e se movimentar ao mesmo tempo.</w:t></w:r></w:p>
'@
$restrRange.InsertXML($restrXml)
